# Generate Report for Archive
#
# The localization status for the zh-cn / de-de targets moved on from
# "Ready for handoff" to "In Translation": update the shared status text
# on the Overview sheet (columns E/F, row 2) and on each locale sheet's
# "Status" column (column C, row 2), then re-fit the now-narrower Status
# columns.

$wb = $excel.ActiveWorkbook

$OLD_STATUS = "Ready for handoff"
$NEW_STATUS = "In Translation"
# `ColumnWidth` (character units) is rounded by Excel to the nearest pixel
# of the workbook's Normal-style font before being stored; 12.5 is the
# input that lands closest to the target stored width (~13.41 chars).
$NEW_STATUS_COLWIDTH = 12.5

# --- 1. Update the status cells -------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("E2").Value = $NEW_STATUS
$overview.Range("F2").Value = $NEW_STATUS
$zhcn.Range("C2").Value = $NEW_STATUS
$dede.Range("C2").Value = $NEW_STATUS

# Safety net: sweep every sheet's used range in case the status shows up
# anywhere else too, so nothing stale is left behind.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # Put the literal on the left: PowerShell's `-eq` coerces the
            # right-hand side to the left-hand side's type, so
            # "$cell.Value2 -eq $OLD_STATUS" would silently coerce our
            # string to a bool (and match TRUE/FALSE cells) whenever
            # Value2 itself is a boolean.
            if ($OLD_STATUS -eq $cell.Value2) {
                $cell.Value = $NEW_STATUS
            }
        }
    }
}

# --- 2. Re-fit the Status columns for the shorter text ---------------------
$overview.Range("E1:F1").ColumnWidth = $NEW_STATUS_COLWIDTH
$zhcn.Range("C1").ColumnWidth = $NEW_STATUS_COLWIDTH
$dede.Range("C1").ColumnWidth = $NEW_STATUS_COLWIDTH
